$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Adam10"
$row2[0,2] = "Epha3"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 52.471316
$row2[0,7] = 157.413948
$row2[0,8] = 0.3942020145328803
$row2[0,9] = 0.3942020145328803
$row2[0,10] = 2
$row2[0,11] = 0.6666666666666666
$row2[0,12] = 0.07579599999999999
$row2[0,13] = 0.227388
$row2[0,14] = 0.001780200955210419
$row2[0,15] = 0.001780200955210419
$row2[0,16] = 3.977115867535999
$row2[0,17] = 35.794042807824
$row2[0,18] = 0.0007017588028173047
$row2[0,19] = 0.0007017588028173048
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Adam10"
$row3[0,2] = "Epha3"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 52.471316
$row3[0,7] = 157.413948
$row3[0,8] = 0.3942020145328803
$row3[0,9] = 0.3942020145328803
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 40.12734033333333
$row3[0,13] = 120.382021
$row3[0,14] = 0.9424604146848589
$row3[0,15] = 0.9424604146848587
$row3[0,16] = 2105.534354869878
$row3[0,17] = 18949.80919382891
$row3[0,18] = 0.3715197940862651
$row3[0,19] = 0.3715197940862651
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Adam10"
$row4[0,2] = "Epha3"
$row4[0,3] = "MuSCs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 52.471316
$row4[0,7] = 157.413948
$row4[0,8] = 0.3942020145328803
$row4[0,9] = 0.3942020145328803
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 2.344072666666666
$row4[0,13] = 7.032217999999999
$row4[0,14] = 0.05505462557763778
$row4[0,15] = 0.05505462557763778
$row4[0,16] = 122.9965776196293
$row4[0,17] = 1106.969198576664
$row4[0,18] = 0.02170264431205825
$row4[0,19] = 0.02170264431205825
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Adam10"
$row5[0,2] = "Epha3"
$row5[0,3] = "Resolving-Mac"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 52.471316
$row5[0,7] = 157.413948
$row5[0,8] = 0.3942020145328803
$row5[0,9] = 0.3942020145328803
$row5[0,10] = 2
$row5[0,11] = 0.6666666666666666
$row5[0,12] = 0.03000666666666667
$row5[0,13] = 0.09002
$row5[0,14] = 0.0007047587822930054
$row5[0,15] = 0.0007047587822930053
$row5[0,16] = 1.574489288773333
$row5[0,17] = 14.17040359896
$row5[0,18] = 0.0002778173317396423
$row5[0,19] = 0.0002778173317396423
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Adam10"
$row6[0,2] = "Epha3"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 23.289271
$row6[0,7] = 69.867813
$row6[0,8] = 0.174965643042042
$row6[0,9] = 0.174965643042042
$row6[0,10] = 2
$row6[0,11] = 0.6666666666666666
$row6[0,12] = 0.07579599999999999
$row6[0,13] = 0.227388
$row6[0,14] = 0.001780200955210419
$row6[0,15] = 0.001780200955210419
$row6[0,16] = 1.765233584716
$row6[0,17] = 15.887102262444
$row6[0,18] = 0.0003114740048724483
$row6[0,19] = 0.0003114740048724484
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Adam10"
$row7[0,2] = "Epha3"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 23.289271
$row7[0,7] = 69.867813
$row7[0,8] = 0.174965643042042
$row7[0,9] = 0.174965643042042
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 40.12734033333333
$row7[0,13] = 120.382021
$row7[0,14] = 0.9424604146848589
$row7[0,15] = 0.9424604146848587
$row7[0,16] = 934.5365035322302
$row7[0,17] = 8410.828531790072
$row7[0,18] = 0.1648981924970059
$row7[0,19] = 0.1648981924970059
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Adam10"
$row8[0,2] = "Epha3"
$row8[0,3] = "MuSCs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 23.289271
$row8[0,7] = 69.867813
$row8[0,8] = 0.174965643042042
$row8[0,9] = 0.174965643042042
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 2.344072666666666
$row8[0,13] = 7.032217999999999
$row8[0,14] = 0.05505462557763778
$row8[0,15] = 0.05505462557763778
$row8[0,16] = 54.59174357769266
$row8[0,17] = 491.325692199234
$row8[0,18] = 0.009632667966630247
$row8[0,19] = 0.009632667966630249
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Adam10"
$row9[0,2] = "Epha3"
$row9[0,3] = "Resolving-Mac"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 23.289271
$row9[0,7] = 69.867813
$row9[0,8] = 0.174965643042042
$row9[0,9] = 0.174965643042042
$row9[0,10] = 2
$row9[0,11] = 0.6666666666666666
$row9[0,12] = 0.03000666666666667
$row9[0,13] = 0.09002
$row9[0,14] = 0.0007047587822930054
$row9[0,15] = 0.0007047587822930053
$row9[0,16] = 0.6988333918066667
$row9[0,17] = 6.28950052626
$row9[0,18] = 0.0001233085735334222
$row9[0,19] = 0.0001233085735334222
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "MuSCs"
$row10[0,1] = "Adam10"
$row10[0,2] = "Epha3"
$row10[0,3] = "ECs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 23.70841633333333
$row10[0,7] = 71.125249
$row10[0,8] = 0.1781145622492915
$row10[0,9] = 0.1781145622492915
$row10[0,10] = 2
$row10[0,11] = 0.6666666666666666
$row10[0,12] = 0.07579599999999999
$row10[0,13] = 0.227388
$row10[0,14] = 0.001780200955210419
$row10[0,15] = 0.001780200955210419
$row10[0,16] = 1.797003124401333
$row10[0,17] = 16.173028119612
$row10[0,18] = 0.0003170797138530742
$row10[0,19] = 0.0003170797138530743
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "MuSCs"
$row11[0,1] = "Adam10"
$row11[0,2] = "Epha3"
$row11[0,3] = "FAPs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 23.70841633333333
$row11[0,7] = 71.125249
$row11[0,8] = 0.1781145622492915
$row11[0,9] = 0.1781145622492915
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 40.12734033333333
$row11[0,13] = 120.382021
$row11[0,14] = 0.9424604146848589
$row11[0,15] = 0.9424604146848587
$row11[0,16] = 951.3556909720253
$row11[0,17] = 8562.201218748229
$row11[0,18] = 0.1678659241988793
$row11[0,19] = 0.1678659241988793
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "MuSCs"
$row12[0,1] = "Adam10"
$row12[0,2] = "Epha3"
$row12[0,3] = "MuSCs"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 23.70841633333333
$row12[0,7] = 71.125249
$row12[0,8] = 0.1781145622492915
$row12[0,9] = 0.1781145622492915
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 2.344072666666666
$row12[0,13] = 7.032217999999999
$row12[0,14] = 0.05505462557763778
$row12[0,15] = 0.05505462557763778
$row12[0,16] = 55.57425069692022
$row12[0,17] = 500.168256272282
$row12[0,18] = 0.009806030534559598
$row12[0,19] = 0.0098060305345596
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "MuSCs"
$row13[0,1] = "Adam10"
$row13[0,2] = "Epha3"
$row13[0,3] = "Resolving-Mac"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 23.70841633333333
$row13[0,7] = 71.125249
$row13[0,8] = 0.1781145622492915
$row13[0,9] = 0.1781145622492915
$row13[0,10] = 2
$row13[0,11] = 0.6666666666666666
$row13[0,12] = 0.03000666666666667
$row13[0,13] = 0.09002
$row13[0,14] = 0.0007047587822930054
$row13[0,15] = 0.0007047587822930053
$row13[0,16] = 0.7114105461088889
$row13[0,17] = 6.40269491498
$row13[0,18] = 0.0001255278019994624
$row13[0,19] = 0.0001255278019994624
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,20
$row14[0,0] = "Resolving-Mac"
$row14[0,1] = "Adam10"
$row14[0,2] = "Epha3"
$row14[0,3] = "ECs"
$row14[0,4] = 3
$row14[0,5] = 1
$row14[0,6] = 33.63867766666667
$row14[0,7] = 100.916033
$row14[0,8] = 0.2527177801757861
$row14[0,9] = 0.2527177801757862
$row14[0,10] = 2
$row14[0,11] = 0.6666666666666666
$row14[0,12] = 0.07579599999999999
$row14[0,13] = 0.227388
$row14[0,14] = 0.001780200955210419
$row14[0,15] = 0.001780200955210419
$row14[0,16] = 2.549677212422666
$row14[0,17] = 22.947094911804
$row14[0,18] = 0.000449888433667591
$row14[0,19] = 0.0004498884336675911
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,20
$row15[0,0] = "Resolving-Mac"
$row15[0,1] = "Adam10"
$row15[0,2] = "Epha3"
$row15[0,3] = "FAPs"
$row15[0,4] = 3
$row15[0,5] = 1
$row15[0,6] = 33.63867766666667
$row15[0,7] = 100.916033
$row15[0,8] = 0.2527177801757861
$row15[0,9] = 0.2527177801757862
$row15[0,10] = 3
$row15[0,11] = 1
$row15[0,12] = 40.12734033333333
$row15[0,13] = 120.382021
$row15[0,14] = 0.9424604146848589
$row15[0,15] = 0.9424604146848587
$row15[0,16] = 1349.830667093632
$row15[0,17] = 12148.47600384269
$row15[0,18] = 0.2381765039027084
$row15[0,19] = 0.2381765039027084
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,20
$row16[0,0] = "Resolving-Mac"
$row16[0,1] = "Adam10"
$row16[0,2] = "Epha3"
$row16[0,3] = "MuSCs"
$row16[0,4] = 3
$row16[0,5] = 1
$row16[0,6] = 33.63867766666667
$row16[0,7] = 100.916033
$row16[0,8] = 0.2527177801757861
$row16[0,9] = 0.2527177801757862
$row16[0,10] = 3
$row16[0,11] = 1
$row16[0,12] = 2.344072666666666
$row16[0,13] = 7.032217999999999
$row16[0,14] = 0.05505462557763778
$row16[0,15] = 0.05505462557763778
$row16[0,16] = 78.85150486124377
$row16[0,17] = 709.663543751194
$row16[0,18] = 0.01391328276438968
$row16[0,19] = 0.01391328276438968
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,20
$row17[0,0] = "Resolving-Mac"
$row17[0,1] = "Adam10"
$row17[0,2] = "Epha3"
$row17[0,3] = "Resolving-Mac"
$row17[0,4] = 3
$row17[0,5] = 1
$row17[0,6] = 33.63867766666667
$row17[0,7] = 100.916033
$row17[0,8] = 0.2527177801757861
$row17[0,9] = 0.2527177801757862
$row17[0,10] = 2
$row17[0,11] = 0.6666666666666666
$row17[0,12] = 0.03000666666666667
$row17[0,13] = 0.09002
$row17[0,14] = 0.0007047587822930054
$row17[0,15] = 0.0007047587822930053
$row17[0,16] = 1.009384587851111
$row17[0,17] = 9.08446129066
$row17[0,18] = 0.0001781050750204785
$row17[0,19] = 0.0001781050750204785
$ws.Range("A17:T17").Value = $row17
